# chore: update Sheets via scheduled runner
#
# Refreshes the cached market-board figures (currentAveragePrice*, Leve
# price/profit columns) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR
# worksheets for the rows whose underlying item data changed.
#
# Each entry: Sheet name, cell reference, new value ($null clears the cell
# entirely, matching rows where a stale profit figure no longer applies).
$edits = @(
    # ALC
    @("ALC", "H15", 32.98),
    @("ALC", "I15", 32.98),
    @("ALC", "K15", 98.94),
    @("ALC", "M15", 70.06),
    @("ALC", "H20", 58870.855),
    @("ALC", "J20", 100024),
    @("ALC", "L20", 100024),
    @("ALC", "N20", -100484),
    @("ALC", "H35", 58870.855),
    @("ALC", "J35", 100024),
    @("ALC", "L35", 100024),
    @("ALC", "N35", -100782),
    @("ALC", "H106", 1971.2858),
    @("ALC", "I106", 1971.2858),
    @("ALC", "J106", 0),
    @("ALC", "K106", 1971.2858),
    @("ALC", "L106", 0),
    @("ALC", "M106", -1340.2858),
    @("ALC", "N106", $null),
    @("ALC", "H129", 842765.5600000001),
    @("ALC", "I129", 419.8),
    @("ALC", "J129", 1090514.2),
    @("ALC", "K129", 1259.4),
    @("ALC", "L129", 3271542.6),
    @("ALC", "M129", 3740.6),
    @("ALC", "N129", -3281542.6),
    @("ALC", "H132", 2042891.4),
    @("ALC", "I132", 2036.1316),
    @("ALC", "J132", 9093119),
    @("ALC", "K132", 6108.3948),
    @("ALC", "L132", 27279357),
    @("ALC", "M132", -3578.3948),
    @("ALC", "N132", -27284417),
    @("ALC", "H137", 1378.3),
    @("ALC", "I137", 996.2),
    @("ALC", "J137", 1760.4),
    @("ALC", "K137", 2988.6),
    @("ALC", "L137", 5281.200000000001),
    @("ALC", "M137", -438.6000000000004),
    @("ALC", "N137", -10381.2),
    @("ALC", "H138", 3265.45),
    @("ALC", "I138", 1646.091),
    @("ALC", "J138", 3722.1924),
    @("ALC", "K138", 4938.272999999999),
    @("ALC", "L138", 11166.5772),
    @("ALC", "M138", 201.7270000000008),
    @("ALC", "N138", -21446.5772),
    # ARM
    @("ARM", "H32", 17703.74),
    @("ARM", "I32", 13013.37),
    @("ARM", "J32", 60503.375),
    @("ARM", "K32", 13013.37),
    @("ARM", "L32", 60503.375),
    @("ARM", "M32", -12726.37),
    @("ARM", "N32", -61077.375),
    @("ARM", "H132", 1773.9183),
    @("ARM", "I132", 1362.3823),
    @("ARM", "J132", 2706.7334),
    @("ARM", "K132", 4087.1469),
    @("ARM", "L132", 8120.2002),
    @("ARM", "M132", -1557.1469),
    @("ARM", "N132", -13180.2002),
    # BSM
    @("BSM", "H141", 54773.5),
    @("BSM", "J141", 54773.5),
    @("BSM", "L141", 54773.5),
    @("BSM", "N141", -65133.5),
    # CRP
    @("CRP", "H31", 2560.2827),
    @("CRP", "I31", 2101.111),
    @("CRP", "J31", 4213.3),
    @("CRP", "K31", 2101.111),
    @("CRP", "L31", 4213.3),
    @("CRP", "M31", -1806.111),
    @("CRP", "N31", -4803.3),
    @("CRP", "H34", 2560.2827),
    @("CRP", "I34", 2101.111),
    @("CRP", "J34", 4213.3),
    @("CRP", "K34", 2101.111),
    @("CRP", "L34", 4213.3),
    @("CRP", "M34", -1899.111),
    @("CRP", "N34", -4617.3),
    @("CRP", "H58", 2898.7144),
    @("CRP", "I58", 725.4358999999999),
    @("CRP", "J58", 11374.5),
    @("CRP", "K58", 725.4358999999999),
    @("CRP", "L58", 11374.5),
    @("CRP", "M58", -522.4358999999999),
    @("CRP", "N58", -11780.5),
    @("CRP", "H127", 40657.145),
    @("CRP", "J127", 40657.145),
    @("CRP", "L127", 40657.145),
    @("CRP", "N127", -50577.145),
    @("CRP", "H136", 2898.7144),
    @("CRP", "I136", 725.4358999999999),
    @("CRP", "J136", 11374.5),
    @("CRP", "K136", 2176.3077),
    @("CRP", "L136", 34123.5),
    @("CRP", "M136", 373.6923000000002),
    @("CRP", "N136", -39223.5),
    # CUL
    @("CUL", "H75", 0),
    @("CUL", "I75", 0),
    @("CUL", "J75", 0),
    @("CUL", "K75", 0),
    @("CUL", "L75", 0),
    @("CUL", "M75", $null),
    @("CUL", "N75", $null),
    @("CUL", "H78", 0),
    @("CUL", "I78", 0),
    @("CUL", "J78", 0),
    @("CUL", "K78", 0),
    @("CUL", "L78", 0),
    @("CUL", "M78", $null),
    @("CUL", "N78", $null),
    @("CUL", "H80", 9279.799999999999),
    @("CUL", "I80", 8199.5),
    @("CUL", "J80", 10000),
    @("CUL", "K80", 24598.5),
    @("CUL", "L80", 30000),
    @("CUL", "M80", -23662.5),
    @("CUL", "N80", -31872),
    @("CUL", "H82", 12820.5),
    @("CUL", "J82", 14022.777),
    @("CUL", "L82", 42068.331),
    @("CUL", "N82", -42880.331),
    @("CUL", "H83", 9279.799999999999),
    @("CUL", "I83", 8199.5),
    @("CUL", "J83", 10000),
    @("CUL", "K83", 73795.5),
    @("CUL", "L83", 90000),
    @("CUL", "M83", -69115.5),
    @("CUL", "N83", -99360),
    @("CUL", "H85", 12820.5),
    @("CUL", "J85", 14022.777),
    @("CUL", "L85", 42068.331),
    @("CUL", "N85", -44876.331),
    @("CUL", "H107", 433380.78),
    @("CUL", "I107", 1764.1428),
    @("CUL", "J107", 708045.9399999999),
    @("CUL", "K107", 5292.428400000001),
    @("CUL", "L107", 2124137.82),
    @("CUL", "M107", -3372.428400000001),
    @("CUL", "N107", -2127977.82),
    @("CUL", "H122", 24494.047),
    @("CUL", "I122", 538.5),
    @("CUL", "J122", 28378.73),
    @("CUL", "K122", 4846.5),
    @("CUL", "L122", 255408.57),
    @("CUL", "M122", -2396.5),
    @("CUL", "N122", -260308.57),
    @("CUL", "H129", 34171.406),
    @("CUL", "I129", 1315),
    @("CUL", "J129", 49106.137),
    @("CUL", "K129", 3945),
    @("CUL", "L129", 147318.411),
    @("CUL", "M129", 1055),
    @("CUL", "N129", -157318.411),
    @("CUL", "H137", 93923.63),
    @("CUL", "I137", 3316),
    @("CUL", "K137", 9948),
    @("CUL", "M137", -4848),
    # GSM
    @("GSM", "H70", 4669.3335),
    @("GSM", "I70", 3999),
    @("GSM", "K70", 3999),
    @("GSM", "M70", -3729),
    @("GSM", "H73", 4669.3335),
    @("GSM", "I73", 3999),
    @("GSM", "K73", 3999),
    @("GSM", "M73", -3063),
    @("GSM", "H102", 2157.7144),
    @("GSM", "I102", 2157.7144),
    @("GSM", "J102", 0),
    @("GSM", "K102", 2157.7144),
    @("GSM", "L102", 0),
    @("GSM", "M102", -535.7143999999998),
    @("GSM", "N102", $null),
    @("GSM", "H126", 3937.3333),
    @("GSM", "I126", 2757.1667),
    @("GSM", "J126", 8658),
    @("GSM", "K126", 8271.500100000001),
    @("GSM", "L126", 25974),
    @("GSM", "M126", -5801.500100000001),
    @("GSM", "N126", -30914),
    # LTW
    @("LTW", "H7", 2324.1936),
    @("LTW", "I7", 2104.5715),
    @("LTW", "J7", 2785.4),
    @("LTW", "K7", 2104.5715),
    @("LTW", "L7", 2785.4),
    @("LTW", "M7", -1992.5715),
    @("LTW", "N7", -3009.4),
    @("LTW", "H126", 2324.1936),
    @("LTW", "I126", 2104.5715),
    @("LTW", "J126", 2785.4),
    @("LTW", "K126", 6313.7145),
    @("LTW", "L126", 8356.200000000001),
    @("LTW", "M126", -3843.7145),
    @("LTW", "N126", -13296.2),
    # WVR
    @("WVR", "H135", 77423.89),
    @("WVR", "J135", 77423.89),
    @("WVR", "L135", 77423.89),
    @("WVR", "N135", -87563.89),
)
$wb = $excel.ActiveWorkbook

# Cache worksheet handles so repeated edits to the same sheet don't
# re-resolve Worksheets.Item() each time.
$sheetCache = @{}

foreach ($edit in $edits) {
    $sheetName = $edit[0]
    $cellRef   = $edit[1]
    $value     = $edit[2]

    if (-not $sheetCache.ContainsKey($sheetName)) {
        $sheetCache[$sheetName] = $wb.Worksheets.Item($sheetName)
    }
    $ws = $sheetCache[$sheetName]
    $ws.Range($cellRef).Value = $value
}
